$wb = $excel.ActiveWorkbook

# --- "Edit Repayment Schedule" sheet: move the remembered selection to B6 ---
# (Selecting a range on a non-active sheet activates that sheet momentarily;
# we re-activate "Repayment schedule" afterwards so the final active tab is correct.)
$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")
[void]$wsEdit.Range("B6").Select()

# --- "Repayment schedule" sheet: add a new (blank) "Variable Instalments" column ---
# Insert a new column before the existing "Late" column (N), shifting
# Late / heading / Outstanding one column to the right (N->O, O->P, P->Q).
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
[void]$wsSchedule.Columns("N").Insert()
$wsSchedule.Columns("N").ColumnWidth = 10.17

# Make "Repayment schedule" the active sheet/tab and set its remembered selection.
[void]$wsSchedule.Activate()
[void]$wsSchedule.Range("K19").Select()
